$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.122.21"
$ws.Range("E2").Value = "  +0.13%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.789.08"
$ws.Range("E3").Value = "  -0.19%  "

# Row 4
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.91"
$ws.Range("E5").Value = "  +0.96%  "

# Row 6
$ws.Range("E6").Value = "  -0.83%  "

# Row 7
$ws.Range("E7").Value = "  +0.07%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.30"
$ws.Range("E8").Value = "  -0.54%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.295"
$ws.Range("E9").Value = "  +3.92%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0687"
$ws.Range("E10").Value = "  -2.73%  "

# Row 11
$ws.Range("E11").Value = "  +1.27%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.047.97"
$ws.Range("E12").Value = "  -0.14%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.46"
$ws.Range("E13").Value = "  +5.86%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.792.18"
$ws.Range("E14").Value = "  -0.18%  "

# Row 15
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.097.81"
$ws.Range("E15").Value = "  +0.13%  "

# Row 16
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.622"
$ws.Range("E16").Value = "  -0.24%  "

# Row 17
$ws.Range("E17").Value = "  +0.45%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.99"
$ws.Range("E18").Value = "  +0.08%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.14"
$ws.Range("E19").Value = "  +0.40%  "

# Row 20
$ws.Range("E20").Value = "  -0.78%  "

# Row 21
$ws.Range("E21").Value = "  +0.13%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.94"
$ws.Range("E22").Value = "  +2.54%  "

# Row 23
$ws.Range("E23").Value = "  +0.40%  "

# Row 24
$ws.Range("E24").Value = "  -2.22%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.41"
$ws.Range("E25").Value = "  +1.67%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.16"
$ws.Range("E26").Value = "  +2.29%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.27"
$ws.Range("E27").Value = "  +0.26%  "

# Row 28
$ws.Range("E28").Value = "  +1.19%  "

# Row 29
$ws.Range("E29").Value = "  +0.23%  "

# Row 30
$ws.Range("E30").Value = "  +1.67%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0517"
$ws.Range("E31").Value = "  -0.12%  "

# Row 32
$ws.Range("E32").Value = "  -0.04%  "

# Row 33
$ws.Range("E33").Value = "  +3.82%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.83"
$ws.Range("E34").Value = "  +1.24%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.403.66"
$ws.Range("E35").Value = "  +1.13%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.652"
$ws.Range("E36").Value = "  +0.99%  "

# Row 37
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0189"
$ws.Range("E37").Value = "  +2.31%  "

# Row 38
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.04"
$ws.Range("E38").Value = "  -0.50%  "

# Row 39
$ws.Range("E39").Value = "  +7.96%  "

# Row 40
$ws.Range("B40").Value = "HuobiToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.39"
$ws.Range("E40").Value = "  +1.49%  "

# Row 41
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "80.02"
$ws.Range("E41").Value = "  +1.21%  "

# Row 42
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.921"
$ws.Range("E42").Value = "  +0.77%  "

# Row 43
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.72"
$ws.Range("E43").Value = "  +0.60%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.42"
$ws.Range("E44").Value = "  +11.96%  "

# Row 45
$ws.Range("E45").Value = "  +0.01%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.06"
$ws.Range("E46").Value = "  +3.59%  "

# Row 47
$ws.Range("E47").Value = "  +2.87%  "

# Row 48
$ws.Range("E48").Value = "  +2.60%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "107.04"
$ws.Range("E49").Value = "  -0.23%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.949.50"
$ws.Range("E50").Value = "  -0.03%  "

# Row 51
$ws.Range("E51").Value = "  +0.23%  "
